# "Modify not to register empty row"
#
# - Rename the "id" header column to "_id" on both the "task" sheet and
#   the "user" sheet (the other columns keep their position/content).
# - Remove the trailing empty row on the "task" sheet (row 4 only had a
#   stray value in M4, with every other cell blank) so an empty row is no
#   longer registered.
# - The "user" sheet becomes the active/selected sheet, with A2 selected;
#   A2 is selected on "task" as well.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("task")
$ws2 = $wb.Worksheets.Item("user")

# --- "task" sheet -----------------------------------------------------
# Header rename: id -> _id
$ws1.Range("A1").Value = "_id"

# Drop the dangling, effectively-empty 4th row (only M4 had a value).
$ws1.Rows.Item(4).Delete()

# Leave the cursor on A2.
$ws1.Range("A2").Select()

# --- "user" sheet -------------------------------------------------------
# Header rename: id -> _id
$ws2.Range("A1").Value = "_id"

# Make "user" the active sheet/tab, with A2 selected.
$ws2.Activate()
$ws2.Range("A2").Select()
